# The deck originally ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral" / "Red Violet" colours (used by the slide master)
#   ppt/theme/theme2.xml -> "Office Theme" / "Office" colours (used by the notes master)
#
# The authored edit swaps the contents of those two theme parts wholesale, so the
# slide master ends up using the "Office Theme" colour palette (and the notes
# master ends up with the former "Integral" palette).
#
# The PowerPoint object model only exposes a single editable Theme for this
# presentation (reached through the slide master), so we reproduce the visible
# half of that swap here: repoint every theme colour slot on
# $p.SlideMaster.Theme.ThemeColorScheme to the "Office Theme" palette, in the
# fixed dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order that PowerPoint uses for
# ThemeColorScheme.Colors(1..12).

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
